$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "73.415.36"
$ws.Range("E2").Value = "  +1.13%  "
$ws.Range("D3").Value = "3.974.68"
$ws.Range("E3").Value = "  -1.86%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D5").Value = "'618.29"
$ws.Range("E5").Value = "  +13.34%  "
$ws.Range("D6").Value = "'168.95"
$ws.Range("E6").Value = "  +10.70%  "
$ws.Range("D7").Value = "'0.680"
$ws.Range("E7").Value = "  -2.56%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +0.14%  "
$ws.Range("D10").Value = "'0.186"
$ws.Range("E10").Value = "  +7.58%  "
$ws.Range("D11").Value = "'55.61"
$ws.Range("E11").Value = "  +3.04%  "
$ws.Range("D12").Value = "'0.0000334"
$ws.Range("E12").Value = "  +1.06%  "
$ws.Range("D13").Value = "'11.16"
$ws.Range("E13").Value = "  +1.66%  "
$ws.Range("D14").Value = "4.613.34"
$ws.Range("E14").Value = "  -1.68%  "
$ws.Range("D15").Value = "3.973.67"
$ws.Range("E15").Value = "  -1.91%  "
$ws.Range("E16").Value = "  +2.27%  "
$ws.Range("D17").Value = "'14.04"
$ws.Range("E17").Value = "  -2.75%  "
$ws.Range("D18").Value = "'20.40"
$ws.Range("E18").Value = "  -2.01%  "
$ws.Range("D19").Value = "73.237.58"
$ws.Range("E19").Value = "  +0.94%  "
$ws.Range("E20").Value = "  -1.05%  "
$ws.Range("D21").Value = "'438.99"
$ws.Range("E21").Value = "  -2.37%  "
$ws.Range("E22").Value = "  +13.43%  "
$ws.Range("D23").Value = "'95.93"
$ws.Range("E23").Value = "  -2.32%  "
$ws.Range("E24").Value = "  -5.07%  "
$ws.Range("D25").Value = "'14.20"
$ws.Range("E25").Value = "  -3.93%  "
$ws.Range("E26").Value = "  -4.32%  "
$ws.Range("D27").Value = "'11.02"
$ws.Range("E27").Value = "  -2.27%  "
$ws.Range("E28").Value = "  +0.20%  "
$ws.Range("D29").Value = "'10.51"
$ws.Range("E29").Value = "  -3.77%  "
$ws.Range("D30").Value = "'36.08"
$ws.Range("E30").Value = "  -3.62%  "
$ws.Range("D31").Value = "'7.80"
$ws.Range("E31").Value = "  -1.86%  "
$ws.Range("E32").Value = "  +0.35%  "
$ws.Range("E33").Value = "  +19.37%  "
$ws.Range("E34").Value = "  -4.07%  "
$ws.Range("D35").Value = "'47.92"
$ws.Range("E35").Value = "  -2.16%  "
$ws.Range("D36").Value = "'70.95"
$ws.Range("E36").Value = "  +5.76%  "
$ws.Range("D37").Value = "'647.61"
$ws.Range("E37").Value = "  -5.41%  "
$ws.Range("E38").Value = "  -4.78%  "
$ws.Range("E39").Value = "  +0.12%  "
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("E41").Value = "  -2.51%  "
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").Value = "'3.24"
$ws.Range("E43").Value = "  -5.50%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "'0.0483"
$ws.Range("E44").Value = "  -2.79%  "
$ws.Range("D45").Value = "'10.58"
$ws.Range("E45").Value = "  -5.91%  "
$ws.Range("D46").Value = "'3.09"
$ws.Range("E46").Value = "  +36.58%  "
$ws.Range("E47").Value = "  -2.42%  "
$ws.Range("D48").Value = "'0.000299"
$ws.Range("E48").Value = "  +7.93%  "
$ws.Range("E49").Value = "  +2.74%  "
$ws.Range("E50").Value = "  -4.73%  "
$ws.Range("D51").Value = "2.817.15"
$ws.Range("E51").Value = "  +2.62%  "
